$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking row): Right column B 5 -> 4, Wrong column C -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total row): Right column B 95 -> 76, Wrong column C -2 -> -4
$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -4

# Update the displayed total text "95 / 140" -> "72 / 112"
$ws.Range("E12").Value = "72 / 112"
